$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 should become a real number (4) instead of a text "4"
$ws.Range("B6").Value = 4

# Add new row 7 with annotation data
$ws.Range("A7").Value = "Ruilin"
# B7 must stay text "3" (not auto-converted to a number) with no cell style,
# so mark it Text, type the value, then reset the style back to Normal.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "3"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "无"
$ws.Range("D7").Value = "FBK"
$ws.Range("E7").Value = "MET"
$ws.Range("F7").Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Range("G7").Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Range("H7").Value = "As you suggested, I did run comparison tests and I will present the results here."
